$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation results table (rows 2-8, columns B-E) with new values
# from the finished base simulation run ("terminacion analisis simulacion base")

$ws.Range("B2").Value = 72.71264502649049
$ws.Range("C2").Value = 75.55286455530215
$ws.Range("D2").Value = 65.74666527964574
$ws.Range("E2").Value = 78.10435175080208

$ws.Range("B3").Value = 95.43673968051225
$ws.Range("C3").Value = 94.55776261496382
$ws.Range("D3").Value = 95.41862177802383
$ws.Range("E3").Value = 95.1653796556581

$ws.Range("B4").Value = 99.38697011759064
$ws.Range("C4").Value = 99.29621816665176
$ws.Range("D4").Value = 99.37212612910128
$ws.Range("E4").Value = 99.39730686782819

$ws.Range("B5").Value = 98.90598092794073
$ws.Range("C5").Value = 98.91215142635062
$ws.Range("D5").Value = 98.88695498066944
$ws.Range("E5").Value = 98.86516365955153

$ws.Range("B6").Value = 98.4910344462462
$ws.Range("C6").Value = 98.41127507373085
$ws.Range("D6").Value = 98.42747415882125
$ws.Range("E6").Value = 98.36970908952426

$ws.Range("B7").Value = 97.43989503442732
$ws.Range("C7").Value = 97.44110361208342
$ws.Range("D7").Value = 97.48030213697916
$ws.Range("E7").Value = 97.45596132465828

$ws.Range("B8").Value = 96.07250952638516
$ws.Range("C8").Value = 96.0372942514505
$ws.Range("D8").Value = 96.03624917568071
$ws.Range("E8").Value = 95.99658762842688
